$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Capture the existing D/E values for the rows that move ---------------
$d16 = $ws.Range("D16").Value()
$e16 = $ws.Range("E16").Value()

$d17 = $ws.Range("D17").Value()
$e17 = $ws.Range("E17").Value()

$d19 = $ws.Range("D19").Value()
$e19 = $ws.Range("E19").Value()

$d20 = $ws.Range("D20").Value()
$e20 = $ws.Range("E20").Value()

$d21 = $ws.Range("D21").Value()
$e21 = $ws.Range("E21").Value()

# --- Row 16 & 17: shift the two-column block one column to the right ------
$ws.Range("E16").Value = $d16
$ws.Range("F16").Value = $e16
$ws.Range("D16").ClearContents()

$ws.Range("E17").Value = $d17
$ws.Range("F17").Value = $e17
$ws.Range("D17").ClearContents()

# --- Rows 19-21: shift the old two-column block to G:H, then add the new --
# --- "Magoosh" half-length review columns in D (and E for row 20) ---------
$ws.Range("G19").Value = $d19
$ws.Range("H19").Value = $e19
$ws.Range("E19").ClearContents()

$ws.Range("G20").Value = $d20
$ws.Range("H20").Value = $e20

$ws.Range("G21").Value = $d21
$ws.Range("H21").Value = $e21
$ws.Range("E21").ClearContents()

$ws.Range("D19").Value = "Magoosh"
$ws.Range("D20").Value = "Magoosh Review"
$ws.Range("D21").Value = "?"
$ws.Range("E20").Value = "BB some"

# --- Conditional formatting tracked the old D20:D21 duplicate check -------
# --- move it along with the data to the new G20:G21 location --------------
$oldRange = $ws.Range("D20:D21")
for ($i = 1; $i -le $oldRange.FormatConditions.Count; $i++) {
    $fc = $oldRange.FormatConditions.Item($i)
    $fc.ModifyAppliesToRange($ws.Range("G20:G21")) | Out-Null
}

# --- View state: scrolled down a bit and a different active selection -----
$win = $excel.ActiveWindow
$win.ScrollRow = 11
$win.ScrollColumn = 1
$ws.Range("D20").Select() | Out-Null
